$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.198.59'
$ws.Range('E2').Value = '  +0.30%  '

$ws.Range('D3').Value = '3.520.65'
$ws.Range('E3').Value = '  +0.11%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '595.34'
$ws.Range('E5').Value = '  +0.40%  '

$ws.Range('D6').Value = '174.62'
$ws.Range('E6').Value = '  +3.39%  '

$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').Value = '0.594'
$ws.Range('E8').Value = '  +2.84%  '

$ws.Range('E9').Value = '  +6.89%  '

$ws.Range('E10').Value = '  -0.07%  '

$ws.Range('D11').Value = '0.438'
$ws.Range('E11').Value = '  -0.43%  '

$ws.Range('D12').Value = '4.124.69'
$ws.Range('E12').Value = '  +0.10%  '

$ws.Range('D13').Value = '0.135'
$ws.Range('E13').Value = '  +0.14%  '

$ws.Range('D14').Value = '28.95'
$ws.Range('E14').Value = '  +2.58%  '

$ws.Range('D15').Value = '0.0000181'
$ws.Range('E15').Value = '  +1.41%  '

$ws.Range('D16').Value = '67.152.78'
$ws.Range('E16').Value = '  +0.35%  '

$ws.Range('D17').Value = '3.517.82'
$ws.Range('E17').Value = '  -0.14%  '

$ws.Range('D18').Value = '6.34'
$ws.Range('E18').Value = '  +0.12%  '

$ws.Range('D19').Value = '14.28'
$ws.Range('E19').Value = '  +1.79%  '

$ws.Range('D20').Value = '396.43'
$ws.Range('E20').Value = '  +1.21%  '

$ws.Range('D21').Value = '7.99'
$ws.Range('E21').Value = '  +0.34%  '

$ws.Range('D22').Value = '73.39'
$ws.Range('E22').Value = '  -0.33%  '

$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Value = '0.542'
$ws.Range('E23').Value = '  +1.67%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.06%  '

$ws.Range('E25').Value = '  -3.46%  '

$ws.Range('D26').Value = '10.20'
$ws.Range('E26').Value = '  +0.10%  '

$ws.Range('D27').Value = '0.180'
$ws.Range('E27').Value = '  -0.70%  '

$ws.Range('D28').Value = '0.997'
$ws.Range('E28').Value = '  -0.77%  '

$ws.Range('E29').Value = '  -2.39%  '

$ws.Range('E30').Value = '  -1.19%  '

$ws.Range('E31').Value = '  +0.38%  '

$ws.Range('D32').Value = '24.04'
$ws.Range('E32').Value = '  +1.82%  '

$ws.Range('D33').Value = '7.42'
$ws.Range('E33').Value = '  -0.74%  '

$ws.Range('E34').Value = '  +3.47%  '

$ws.Range('D35').Value = '163.73'
$ws.Range('E35').Value = '  +1.59%  '

$ws.Range('D36').Value = '0.896'
$ws.Range('E36').Value = '  -0.55%  '

$ws.Range('D37').Value = '1.92'
$ws.Range('E37').Value = '  -0.72%  '

$ws.Range('D38').Value = '6.92'
$ws.Range('E38').Value = '  +3.13%  '

$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = '27.69'
$ws.Range('E39').Value = '  +4.84%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '4.70'
$ws.Range('E40').Value = '  +0.74%  '

$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '0.0744'
$ws.Range('E41').Value = '  -0.94%  '

$ws.Range('D42').Value = '26.52'
$ws.Range('E42').Value = '  +0.10%  '

$ws.Range('E43').Value = '  +3.23%  '

$ws.Range('D44').Value = '2.802.04'
$ws.Range('E44').Value = '  -1.09%  '

$ws.Range('D45').Value = '42.88'

$ws.Range('D46').Value = '0.0306'
$ws.Range('E46').Value = '  -2.85%  '

$ws.Range('D47').Value = '339.52'
$ws.Range('E47').Value = '  -4.03%  '

$ws.Range('D48').Value = '1.09'
$ws.Range('E48').Value = '  +0.32%  '

$ws.Range('E49').Value = '  -0.34%  '

$ws.Range('D50').Value = '6.53'
$ws.Range('E50').Value = '  +0.07%  '

$ws.Range('D51').Value = '0.847'
$ws.Range('E51').Value = '  -1.26%  '
